# Add season "record" columns (Wins / Losses / Ties) to the team stats sheet.
# The old scraper only pulled team statistics and missed the season record,
# so this populates three new trailing columns (AC:AE) with the team's
# win/loss/tie totals for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells should look exactly like the existing header row (A1:AB1):
# bold font, thin border, centered/top-aligned. Copy the formatting from A1
# so the new cells reuse the same style instead of minting a near-duplicate.
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Every player row (2-48) gets the same 1996 Pittsburgh Pirates season
# record: 73 wins, 89 losses, 0 ties.
$lastRow = 48
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 73
    $ws.Cells.Item($r, 30).Value = 89
    $ws.Cells.Item($r, 31).Value = 0
}
